$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 195.0069102805903
$ws.Range("C2").Value = 16.37240170114682
$ws.Range("D2").Value = 223.6666899063773
$ws.Range("F2").Value = 3046.395281440311
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.02452288890473963
$ws.Range("K2").Value = 0.06221187317580171
$ws.Range("L2").Value = 0.06004362877703256
$ws.Range("N2").Value = 1001.546187016945
$ws.Range("P2").Value = 0.8941876865246061
$ws.Range("Q2").Value = 200
$ws.Range("S2").Value = 0.8296150187497769
$ws.Range("T2").Value = 2527.335278531336
$ws.Range("V2").Value = 0.3672251794272522
$ws.Range("W2").Value = 1001.546187016945
$ws.Range("Y2").Value = 0.4767301622731319
$ws.Range("B3").Value = 114.0607491686329
$ws.Range("C3").Value = 15.55963215570908
$ws.Range("D3").Value = 212.5632808156979
$ws.Range("F3").Value = 2997.150922555896
$ws.Range("J3").Value = 0.00526963314993755
$ws.Range("K3").Value = 0.03768265982457567
$ws.Range("L3").Value = 0.02434972138313655
$ws.Range("N3").Value = 996.9122513350209
$ws.Range("P3").Value = 0.9408962791339733
$ws.Range("Q3").Value = 200
$ws.Range("S3").Value = 0.8157821215605918
$ws.Range("T3").Value = 2445.022138239934
$ws.Range("V3").Value = 0.3769012882434292
$ws.Range("W3").Value = 996.9122513350209
$ws.Range("Y3").Value = 0.4460531942133563
$ws.Range("B4").Value = 141.8694817229864
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 212.352668818627
$ws.Range("F4").Value = 2430.081564951446
$ws.Range("H4").Value = 611.3946160619125
$ws.Range("J4").Value = 0.0257223262103678
$ws.Range("K4").Value = 0.05169881595332022
$ws.Range("L4").Value = 0.01762354008841846
$ws.Range("N4").Value = 974.7074151959887
$ws.Range("P4").Value = 0.9418294628113313
$ws.Range("Q4").Value = 200
$ws.Range("S4").Value = 0.8209452698122482
$ws.Range("T4").Value = 1994.963966004835
$ws.Range("V4").Value = 0.3981169233047824
$ws.Range("W4").Value = 974.7074151959887
$ws.Range("Y4").Value = 0.414348410058319
$ws.Range("Z4").Value = 253.3303870834698
$ws.Range("B5").Value = 150.8817711372159
$ws.Range("C5").Value = 40.00000000000001
$ws.Range("D5").Value = 217.3042912822645
$ws.Range("F5").Value = 2508.232377046967
$ws.Range("H5").Value = 602.333146953456
$ws.Range("J5").Value = 0.004488578146070251
$ws.Range("K5").Value = 0.05323732373615354
$ws.Range("L5").Value = 0.02718562786190657
$ws.Range("N5").Value = 1006.086638669547
$ws.Range("P5").Value = 0.920368386743972
$ws.Range("Q5").Value = 200
$ws.Range("S5").Value = 0.8251549975806217
$ws.Range("T5").Value = 2069.680481013827
$ws.Range("V5").Value = 0.3931873223339269
$ws.Range("W5").Value = 1006.086638669547
$ws.Range("Y5").Value = 0.4799947824479487
$ws.Range("Z5").Value = 289.1167678331124
$ws.Range("B6").Value = 206.0151730223007
$ws.Range("C6").Value = 40.00000000000001
$ws.Range("D6").Value = 217.1579832821567
$ws.Range("F6").Value = 2736.594889309589
$ws.Range("H6").Value = 602.6008905936533
$ws.Range("J6").Value = 0.02448315454606691
$ws.Range("K6").Value = 0.06858679466424583
$ws.Range("L6").Value = 0.02158010219928313
$ws.Range("N6").Value = 1023.739735273847
$ws.Range("P6").Value = 0.9209884756580047
$ws.Range("Q6").Value = 200
$ws.Range("S6").Value = 0.8519370354569534
$ws.Range("T6").Value = 2331.406537245061
$ws.Range("V6").Value = 0.3578591399661504
$ws.Range("W6").Value = 1023.739735273847
$ws.Range("Y6").Value = 0.5465097415581022
$ws.Range("Z6").Value = 329.3272569810197
$ws.Range("B7").Value = 210.3943289132366
$ws.Range("C7").Value = 40
$ws.Range("D7").Value = 210.8931742376801
$ws.Range("F7").Value = 2524.48167866078
$ws.Range("H7").Value = 614.0654911450455
$ws.Range("J7").Value = 0.03039519206980118
$ws.Range("K7").Value = 0.075024350053667
$ws.Range("L7").Value = 0.02375413277980732
$ws.Range("N7").Value = 1020.5778852583
$ws.Range("P7").Value = 0.9483474309822694
$ws.Range("Q7").Value = 200
$ws.Range("S7").Value = 0.8620282240290101
$ws.Range("T7").Value = 2176.174458049726
$ws.Range("V7").Value = 0.3882416450517856
$ws.Range("W7").Value = 1020.5778852583
$ws.Range("Y7").Value = 0.4112651881553938
$ws.Range("Z7").Value = 252.5437597555014
$ws.Range("B8").Value = 135.7438825309016
$ws.Range("C8").Value = 40.00000000000001
$ws.Range("D8").Value = 225.9532802531134
$ws.Range("F8").Value = 2263.920682801419
$ws.Range("H8").Value = 586.5054971368027
$ws.Range("J8").Value = 0.01834521704604285
$ws.Range("K8").Value = 0.05343249263847207
$ws.Range("L8").Value = 0.01812735839037085
$ws.Range("N8").Value = 1012.381961101153
$ws.Range("P8").Value = 0.8851387321129374
$ws.Range("Q8").Value = 200
$ws.Range("S8").Value = 0.8470047133927943
$ws.Range("T8").Value = 1917.551489080235
$ws.Range("V8").Value = 0.4166699619360318
$ws.Range("W8").Value = 1012.381961101153
$ws.Range("Y8").Value = 0.5322132588149023
$ws.Range("Z8").Value = 312.1460019440321
$ws.Range("B9").Value = 112.0608846513275
$ws.Range("C9").Value = 40
$ws.Range("F9").Value = 2064.293057906641
$ws.Range("H9").Value = 1000
$ws.Range("J9").Value = 0.02297127951841579
$ws.Range("K9").Value = 0.05276057339089033
$ws.Range("L9").Value = 0.003147599269338761
$ws.Range("N9").Value = 1016.571325306001
$ws.Range("P9").Value = 0.9015562974029467
$ws.Range("S9").Value = 0.8578027257244812
$ws.Range("T9").Value = 1770.756211766441
$ws.Range("V9").Value = 0.4354257695562926
$ws.Range("W9").Value = 1016.571325306001
$ws.Range("Y9").Value = 0.5639042433138978
$ws.Range("Z9").Value = 563.9042433138977
$ws.Range("B10").Value = 131.3777617639103
$ws.Range("C10").Value = 16.33907086210717
$ws.Range("D10").Value = 223.2113505752345
$ws.Range("F10").Value = 2601.679498268624
$ws.Range("J10").Value = 0.02856099409530397
$ws.Range("K10").Value = 0.04804689577686777
$ws.Range("L10").Value = 0.03233245612175138
$ws.Range("N10").Value = 1002.420838364033
$ws.Range("P10").Value = 0.8960117820378897
$ws.Range("Q10").Value = 200
$ws.Range("S10").Value = 0.8420504339628174
$ws.Range("T10").Value = 2190.74535054926
$ws.Range("V10").Value = 0.4192921835584589
$ws.Range("W10").Value = 1002.420838364033
$ws.Range("Y10").Value = 0.4832423747172801
$ws.Range("B11").Value = 157.207514328091
$ws.Range("C11").Value = 40
$ws.Range("F11").Value = 2114.213088858497
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 0.02625684770907999
$ws.Range("K11").Value = 0.06612645055040436
$ws.Range("L11").Value = 0.01740210705467193
$ws.Range("N11").Value = 1030.005553117358
$ws.Range("P11").Value = 0.9388372227906356
$ws.Range("S11").Value = 0.8475389390374807
$ws.Range("T11").Value = 1791.877918230285
$ws.Range("V11").Value = 0.4375480352361826
$ws.Range("W11").Value = 1030.005553117358
$ws.Range("Y11").Value = 0.5621620274899031
$ws.Range("Z11").Value = 562.162027489903
